$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.795.98'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.386.34'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.23'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.65'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.387.44'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.59'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('E11').Value = '  -3.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.387'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.957.24'
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.41'
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000173'
$ws.Range('E16').Value = '  -3.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.384.19'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.962.57'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.04'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.84'
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.39'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.67'
$ws.Range('E22').Value = '  -3.44%  '
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.512.03'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000125'
$ws.Range('E26').Value = '  -2.75%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.16'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.78'
$ws.Range('E28').Value = '  +11.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.54'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.79'
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.22'
$ws.Range('E36').Value = '  -5.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.83'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.80'
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0756'
$ws.Range('E40').Value = '  -4.57%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.773'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.21'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.70'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.23'
$ws.Range('E47').Value = '  -4.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.38'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.81'
$ws.Range('E49').Value = '  -2.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.385.78'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.42'
$ws.Range('E51').Value = '  +4.77%  '
